$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new blank column before the
# existing "Late" column (old column N / 14th column), shifting
# Late / Date / Outstanding one column to the right. ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert()

# Carry over the width of the column immediately to the left (M) onto
# the freshly inserted column, matching Excel's own insert behaviour.
$ws.Columns.Item(14).ColumnWidth = 9.83

# Move the selection to a single cell, as in the edited workbook.
$ws.Range("R8").Select() | Out-Null

# Make "Repayment schedule" the active sheet/tab (this also clears the
# tabSelected flag that was previously on "Acc_Disbursement").
$ws.Activate() | Out-Null
